# Update the "想去人数" (want-to-go count) figures in the F column for both
# the "展览" sheet and the "全部类型" sheet, reflecting the regenerated
# gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" - F-column row numbers as they appear on that sheet.
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 105
$ws1.Range("F13").Value = 2486
$ws1.Range("F14").Value = 56
$ws1.Range("F20").Value = 586
$ws1.Range("F21").Value = 176
$ws1.Range("F25").Value = 2081
$ws1.Range("F26").Value = 4181
$ws1.Range("F30").Value = 1225
$ws1.Range("F32").Value = 2127
$ws1.Range("F39").Value = 722
$ws1.Range("F42").Value = 6
$ws1.Range("F43").Value = 432

# Sheet "全部类型" - same events, but row numbers are shifted by one for
# several rows because this sheet carries one extra row versus "展览".
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 105
$ws4.Range("F13").Value = 2486
$ws4.Range("F14").Value = 56
$ws4.Range("F21").Value = 586
$ws4.Range("F22").Value = 176
$ws4.Range("F26").Value = 2081
$ws4.Range("F27").Value = 4181
$ws4.Range("F31").Value = 1225
$ws4.Range("F33").Value = 2127
$ws4.Range("F40").Value = 722
$ws4.Range("F43").Value = 6
$ws4.Range("F44").Value = 432
